$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update text "NEBNextPoly(A)E7490" -> "NEBNextPoly(A)E7490L" (column G, rows 2-19)
for ($r = 2; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "NEBNextPoly(A)E7490") {
        $cell.Value = "NEBNextPoly(A)E7490L"
    }
}

# 2. Widen column G (7) to fit the longer text, leaving the other columns untouched
$ws.Columns.Item(7).ColumnWidth = 42.14

# 3. Normalize the font color on G3:G19 to match G2's explicit-black styling
for ($r = 3; $r -le 19; $r++) {
    $ws.Cells.Item($r, 7).Font.Color = 0
}

# 4. Rewrite column I (rows 2-19) so the constant boolean is backed by a formula
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 9).Formula = "=FALSE()"
}

# 5. Move the active selection from I2:I19 to G2:G19
$ws.Range("G2:G19").Select()
